# Doing Updates for Financials
# Updates a set of historical figures on the PDPTF financials sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDPTF")

# Row 14 - Inventory
$ws.Range("I14").Value = 3900

# Row 17 - Total Current Assets
$ws.Range("D17").Value = 2100
$ws.Range("E17").Value = 1600
$ws.Range("F17").Value = 2000
$ws.Range("I17").Value = 4100
$ws.Range("J17").Value = 200

# Row 18 - Long Term Investments
$ws.Range("E18").Value = -900
$ws.Range("F18").Value = -1900

# Row 21 - Other Assets
$ws.Range("F21").Value = -1900

# Row 23 - Total Assets
$ws.Range("D23").Value = -1600
$ws.Range("F23").Value = -1900
$ws.Range("I23").Value = -4100

# Row 26 - Total Current Liabilities
$ws.Range("D26").Value = -1600
$ws.Range("F26").Value = -1900
$ws.Range("I26").Value = -3900

# Row 27 - Long Term Debt
$ws.Range("D27").Value = -1600
$ws.Range("F27").Value = -1900
$ws.Range("I27").Value = -3900

# Row 33 - Net Income
$ws.Range("D33").Value = -1600
$ws.Range("F33").Value = -1900
$ws.Range("I33").Value = -3900

# Row 35 - Net Income Applicable To Common Shares
$ws.Range("D35").Value = -1600
$ws.Range("F35").Value = -1900
$ws.Range("I35").Value = -3900

# Row 48 - Total Assets (Balance Sheet section)
$ws.Range("J48").Value = 3900

# Row 54 - Total Liabilities
$ws.Range("J54").Value = 4100

# Row 57 - Accounts Payable
$ws.Range("H57").Value = 100

# Row 60 - Total Current Liabilities
$ws.Range("D60").Value = 200

# Row 66 - Total Liabilities
$ws.Range("E66").Value = 400

# Row 72 - Retained Earnings
$ws.Range("D72").Value = -3900
$ws.Range("E72").Value = -2700
$ws.Range("J72").Value = -8300

# Row 76 - Total Stockholder Equity
$ws.Range("J76").Value = 4000

# Row 81 - Net Income (Cash Flow Statement section)
$ws.Range("D81").Value = -1600
$ws.Range("F81").Value = -1900
$ws.Range("I81").Value = -3900

# Row 89 - Total Cash Flow From Operating Activities
$ws.Range("D89").Value = -800

# Row 100 - Net Borrowings
$ws.Range("F100").Value = 1400
